$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold numeric-looking text that Excel
# would otherwise auto-convert to a Number/Percentage on assignment, which
# would also silently add a new cell style. Force a Text number format before
# the write, then restore the default "Normal" style so the cell ends up with
# no style index at all, matching the original formatting.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '303.20'
Set-TextValue "E2" '-4.97%'
Set-TextValue "D3" '35.04'
Set-TextValue "E3" '-2.99%'
Set-TextValue "D4" '5.066'
Set-TextValue "E4" '-1.26%'
Set-TextValue "D5" '0.07983'
Set-TextValue "E5" '-2.75%'
Set-TextValue "D6" '1.929'
Set-TextValue "E6" '-10.14%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue "D7" '7.748'
Set-TextValue "E7" '-3.23%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D8" '2.942'
Set-TextValue "E8" '5.04%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D9" '0.9216'
Set-TextValue "E9" '-0.62%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D10" '0.1222'
Set-TextValue "E10" '22.48%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D11" '0.1844'
Set-TextValue "E11" '-2.37%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D12" '0.09357'
Set-TextValue "E12" '1.60%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D13" '0.03526'
Set-TextValue "E13" '-2.50%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D14" '0.09854'
Set-TextValue "E14" '-0.63%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D15" '0.001386'
Set-TextValue "E15" '-3.56%'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue "D16" '0.04503'
Set-TextValue "E16" '-2.02%'
Set-TextValue "D17" '0.005862'
Set-TextValue "E17" '3.47%'
Set-TextValue "D18" '3.496'
Set-TextValue "E18" '0.82%'
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D19" '4.056'
Set-TextValue "E19" '-1.92%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue "D20" '0.3447'
Set-TextValue "E20" '2.10%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue "D21" '0.1291'
Set-TextValue "E21" '-0.78%'
$ws.Range("B22").Value = 'MCDex'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue "D22" '5.039'
Set-TextValue "E22" '-0.37%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue "D23" '0.2464'
Set-TextValue "E23" '12.55%'
Set-TextValue "D24" '0.001217'
Set-TextValue "E24" '-2.41%'
Set-TextValue "D25" '0.004856'
Set-TextValue "E25" '2.57%'
Set-TextValue "E26" '-0.18%'
Set-TextValue "E27" '-6.86%'
Set-TextValue "D39" '0.01927'
Set-TextValue "E39" '-3.81%'
Set-TextValue "D40" '0.04751'
Set-TextValue "E40" '-4.56%'
Set-TextValue "D41" '0.007519'
Set-TextValue "E41" '-3.18%'
Set-TextValue "D42" '0.009549'
Set-TextValue "E42" '22.21%'
Set-TextValue "D43" '0.1329'
Set-TextValue "E43" '-5.06%'
Set-TextValue "D44" '0.002108'
Set-TextValue "E44" '0.48%'
Set-TextValue "D45" '0.01116'
Set-TextValue "E45" '-5.99%'
Set-TextValue "D46" '0.00006295'
Set-TextValue "E46" '-2.53%'
Set-TextValue "E47" '-0.15%'
Set-TextValue "E48" '120.94%'
Set-TextValue "E49" '-31.38%'
Set-TextValue "E50" '-0.15%'
Set-TextValue "E51" '-0.15%'
